$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 193; this pushes the existing rows
# 193-201 down to 194-202 (matching the weekly-update diff).
$ws.Rows("193").Insert()

# Populate the newly inserted row 193 with this week's record.
$ws.Range("A193").Value = 7
$ws.Range("B193").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C193").Value = "Ñuble"
$ws.Range("D193").Value = 45166
$ws.Range("E193").Value = 16
$ws.Range("F193").Value = 100112037
$ws.Range("G193").Value = "Cebollín"
$ws.Range("H193").Value = "Sin especificar"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 250
$ws.Range("K193").Value = 6000
$ws.Range("L193").Value = 6000
$ws.Range("M193").Value = 6000
$ws.Range("N193").Value = "$/paquete 36 unidades"
$ws.Range("O193").Value = "Provincia de Diguillín"
$ws.Range("P193").Value = 167
$ws.Range("Q193").Value = 36
$ws.Range("R193").Value = "Hortaliza"
